$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.800.50"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "3.281.70"
$ws.Range("E3").Value = "  -1.35%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'573.37"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").Value = "'174.97"
$ws.Range("E6").Value = "  -5.64%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("D9").Value = "3.277.89"
$ws.Range("E9").Value = "  -1.26%  "

# Row 10
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  -3.69%  "

# Row 11
$ws.Range("D11").Value = "'0.572"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12
$ws.Range("D12").Value = "'45.38"
$ws.Range("E12").Value = "  -3.64%  "

# Row 13
$ws.Range("D13").Value = "'0.0000267"
$ws.Range("E13").Value = "  -0.32%  "

# Row 14
$ws.Range("D14").Value = "'684.89"
$ws.Range("E14").Value = "  +3.74%  "

# Row 15
$ws.Range("D15").Value = "3.811.18"
$ws.Range("E15").Value = "  -1.19%  "

# Row 16
$ws.Range("D16").Value = "'8.27"
$ws.Range("E16").Value = "  -2.44%  "

# Row 17
$ws.Range("D17").Value = "66.951.19"
$ws.Range("E17").Value = "  +1.16%  "

# Row 18
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").Value = "3.285.64"
$ws.Range("E19").Value = "  -1.14%  "

# Row 20
$ws.Range("D20").Value = "'17.26"
$ws.Range("E20").Value = "  -3.78%  "

# Row 21
$ws.Range("D21").Value = "'10.70"
$ws.Range("E21").Value = "  -3.28%  "

# Row 22
$ws.Range("D22").Value = "'0.886"
$ws.Range("E22").Value = "  -1.22%  "

# Row 23
$ws.Range("D23").Value = "'16.98"
$ws.Range("E23").Value = "  -5.34%  "

# Row 24
$ws.Range("D24").Value = "'5.14"
$ws.Range("E24").Value = "  +2.24%  "

# Row 25
$ws.Range("D25").Value = "'98.62"
$ws.Range("E25").Value = "  -2.39%  "

# Row 26
$ws.Range("E26").Value = "  -3.04%  "

# Row 27
$ws.Range("E27").Value = "  -3.16%  "

# Row 28
$ws.Range("D28").Value = "'9.25"
$ws.Range("E28").Value = "  -2.85%  "

# Row 29
$ws.Range("D29").Value = "'33.52"
$ws.Range("E29").Value = "  +6.81%  "

# Row 30
$ws.Range("D30").Value = "'8.36"
$ws.Range("E30").Value = "  -1.28%  "

# Row 31
$ws.Range("D31").Value = "'6.72"
$ws.Range("E31").Value = "  +0.58%  "

# Row 32
$ws.Range("D32").Value = "'570.08"
$ws.Range("E32").Value = "  -4.04%  "

# Row 33
$ws.Range("D33").Value = "3.871.28"
$ws.Range("E33").Value = "  +0.80%  "

# Row 34
$ws.Range("D34").Value = "'10.81"
$ws.Range("E34").Value = "  -1.64%  "

# Row 35
$ws.Range("E35").Value = "  -3.07%  "

# Row 37
$ws.Range("D37").Value = "'55.24"
$ws.Range("E37").Value = "  -1.01%  "

# Row 38
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  -14.80%  "

# Row 39
$ws.Range("E39").Value = "  +0.80%  "

# Row 40
$ws.Range("D40").Value = "'3.38"
$ws.Range("E40").Value = "  -1.44%  "

# Row 41
$ws.Range("D41").Value = "'2.58"
$ws.Range("E41").Value = "  -3.63%  "

# Row 42
$ws.Range("D42").Value = "'31.74"
$ws.Range("E42").Value = "  -3.38%  "

# Row 43
$ws.Range("D43").Value = "0.0₃0665"
$ws.Range("E43").Value = "  -4.69%  "

# Row 44
$ws.Range("D44").Value = "'0.325"
$ws.Range("E44").Value = "  -3.18%  "

# Row 45
$ws.Range("D45").Value = "'2.98"
$ws.Range("E45").Value = "  -6.45%  "

# Row 46
$ws.Range("D46").Value = "'0.0402"
$ws.Range("E46").Value = "  -2.13%  "

# Row 47
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.01"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.127"
$ws.Range("E48").Value = "  -0.58%  "

# Row 49
$ws.Range("E49").Value = "  -0.65%  "

# Row 50
$ws.Range("E50").Value = "  +6.03%  "

# Row 51
$ws.Range("D51").Value = "'129.96"
$ws.Range("E51").Value = "  -0.36%  "
